$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): "Name" / "Team Name" ---
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Team Name"

# --- Data row (row 2): "John joe" / "Team 01" ---
$ws.Range("A2").Value = "John joe"
$ws.Range("B2").Value = "Team 01"

# Apply the 12pt font to the data rows (2-3) first, then 16pt to the
# header row, so the shared font table / shared string table end up in
# the same order as the source workbook.
$ws.Range("A2:B3").Font.Size = 12
$ws.Range("A1:B1").Font.Size = 16

# Row heights to match the bigger header font and the 12pt body font.
$ws.Range("1:1").RowHeight = 21
$ws.Range("2:3").RowHeight = 15.75

# Widen column B to fit "Team Name" / "Team 01".
$ws.Columns("B").ColumnWidth = 16.140625

# Leave the selection where Excel would land after typing the last row.
$ws.Range("B4").Select() | Out-Null
